# This sheet contains weekly price data for "Perejil" (Hortaliza) at Feria Lagunitas
# de Puerto Montt. The update shifts the existing series down by two rows (rows
# 176-228 now hold what used to be rows 174-226), inserts two new most-recent weekly
# records at rows 174-175, and the two rows that fall off the end of the shift are
# appended as brand-new rows 227-228, extending the table from R226 to R228.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colD = 4; $colJ = 10; $colK = 11; $colL = 12; $colM = 13
$colN = 14; $colO = 15; $colP = 16; $colQ = 17

$ws.Cells.Item(174, $colD).Value = 44642
$ws.Cells.Item(174, $colJ).Value = 80
$ws.Cells.Item(174, $colK).Value = 5000
$ws.Cells.Item(174, $colL).Value = 5000
$ws.Cells.Item(174, $colM).Value = 5000
$ws.Cells.Item(174, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(174, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(174, $colP).Value = 2500
$ws.Cells.Item(174, $colQ).Value = 2

$ws.Cells.Item(175, $colD).Value = 44642
$ws.Cells.Item(175, $colJ).Value = 180
$ws.Cells.Item(175, $colK).Value = 5000
$ws.Cells.Item(175, $colL).Value = 5000
$ws.Cells.Item(175, $colM).Value = 5000
$ws.Cells.Item(175, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(175, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(175, $colP).Value = 1667
$ws.Cells.Item(175, $colQ).Value = 3

$ws.Cells.Item(176, $colD).Value = 44473
$ws.Cells.Item(176, $colJ).Value = 120
$ws.Cells.Item(176, $colK).Value = 4500
$ws.Cells.Item(176, $colL).Value = 4500
$ws.Cells.Item(176, $colM).Value = 4500
$ws.Cells.Item(176, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(176, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(176, $colP).Value = 1500
$ws.Cells.Item(176, $colQ).Value = 3

$ws.Cells.Item(177, $colD).Value = 44357
$ws.Cells.Item(177, $colJ).Value = 40
$ws.Cells.Item(177, $colK).Value = 3500
$ws.Cells.Item(177, $colL).Value = 3500
$ws.Cells.Item(177, $colM).Value = 3500
$ws.Cells.Item(177, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(177, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(177, $colP).Value = 1167
$ws.Cells.Item(177, $colQ).Value = 3

$ws.Cells.Item(178, $colD).Value = 44537
$ws.Cells.Item(178, $colJ).Value = 180
$ws.Cells.Item(178, $colK).Value = 5000
$ws.Cells.Item(178, $colL).Value = 5000
$ws.Cells.Item(178, $colM).Value = 5000
$ws.Cells.Item(178, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(178, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(178, $colP).Value = 1667
$ws.Cells.Item(178, $colQ).Value = 3

$ws.Cells.Item(179, $colD).Value = 44553
$ws.Cells.Item(179, $colJ).Value = 80
$ws.Cells.Item(179, $colK).Value = 7000
$ws.Cells.Item(179, $colL).Value = 7000
$ws.Cells.Item(179, $colM).Value = 7000
$ws.Cells.Item(179, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(179, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(179, $colP).Value = 3500
$ws.Cells.Item(179, $colQ).Value = 2

$ws.Cells.Item(180, $colD).Value = 44490
$ws.Cells.Item(180, $colJ).Value = 60
$ws.Cells.Item(180, $colK).Value = 5000
$ws.Cells.Item(180, $colL).Value = 5000
$ws.Cells.Item(180, $colM).Value = 5000
$ws.Cells.Item(180, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(180, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(180, $colP).Value = 1667
$ws.Cells.Item(180, $colQ).Value = 3

$ws.Cells.Item(181, $colD).Value = 44397
$ws.Cells.Item(181, $colJ).Value = 160
$ws.Cells.Item(181, $colK).Value = 5000
$ws.Cells.Item(181, $colL).Value = 5000
$ws.Cells.Item(181, $colM).Value = 5000
$ws.Cells.Item(181, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(181, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(181, $colP).Value = 1667
$ws.Cells.Item(181, $colQ).Value = 3

$ws.Cells.Item(182, $colD).Value = 44446
$ws.Cells.Item(182, $colJ).Value = 160
$ws.Cells.Item(182, $colK).Value = 5000
$ws.Cells.Item(182, $colL).Value = 5000
$ws.Cells.Item(182, $colM).Value = 5000
$ws.Cells.Item(182, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(182, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(182, $colP).Value = 1667
$ws.Cells.Item(182, $colQ).Value = 3

$ws.Cells.Item(183, $colD).Value = 44641
$ws.Cells.Item(183, $colJ).Value = 70
$ws.Cells.Item(183, $colK).Value = 5000
$ws.Cells.Item(183, $colL).Value = 5000
$ws.Cells.Item(183, $colM).Value = 5000
$ws.Cells.Item(183, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(183, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(183, $colP).Value = 2500
$ws.Cells.Item(183, $colQ).Value = 2

$ws.Cells.Item(184, $colD).Value = 44421
$ws.Cells.Item(184, $colJ).Value = 180
$ws.Cells.Item(184, $colK).Value = 5000
$ws.Cells.Item(184, $colL).Value = 5000
$ws.Cells.Item(184, $colM).Value = 5000
$ws.Cells.Item(184, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(184, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(184, $colP).Value = 1667
$ws.Cells.Item(184, $colQ).Value = 3

$ws.Cells.Item(185, $colD).Value = 44329
$ws.Cells.Item(185, $colJ).Value = 50
$ws.Cells.Item(185, $colK).Value = 4500
$ws.Cells.Item(185, $colL).Value = 4500
$ws.Cells.Item(185, $colM).Value = 4500
$ws.Cells.Item(185, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(185, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(185, $colP).Value = 1500
$ws.Cells.Item(185, $colQ).Value = 3

$ws.Cells.Item(186, $colD).Value = 44637
$ws.Cells.Item(186, $colJ).Value = 80
$ws.Cells.Item(186, $colK).Value = 5000
$ws.Cells.Item(186, $colL).Value = 5000
$ws.Cells.Item(186, $colM).Value = 5000
$ws.Cells.Item(186, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(186, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(186, $colP).Value = 1667
$ws.Cells.Item(186, $colQ).Value = 3

$ws.Cells.Item(187, $colD).Value = 44208
$ws.Cells.Item(187, $colJ).Value = 100
$ws.Cells.Item(187, $colK).Value = 5000
$ws.Cells.Item(187, $colL).Value = 5000
$ws.Cells.Item(187, $colM).Value = 5000
$ws.Cells.Item(187, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(187, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(187, $colP).Value = 2500
$ws.Cells.Item(187, $colQ).Value = 2

$ws.Cells.Item(188, $colD).Value = 44355
$ws.Cells.Item(188, $colJ).Value = 160
$ws.Cells.Item(188, $colK).Value = 3500
$ws.Cells.Item(188, $colL).Value = 3500
$ws.Cells.Item(188, $colM).Value = 3500
$ws.Cells.Item(188, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(188, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(188, $colP).Value = 1167
$ws.Cells.Item(188, $colQ).Value = 3

$ws.Cells.Item(189, $colD).Value = 44530
$ws.Cells.Item(189, $colJ).Value = 120
$ws.Cells.Item(189, $colK).Value = 6000
$ws.Cells.Item(189, $colL).Value = 6000
$ws.Cells.Item(189, $colM).Value = 6000
$ws.Cells.Item(189, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(189, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(189, $colP).Value = 3000
$ws.Cells.Item(189, $colQ).Value = 2

$ws.Cells.Item(190, $colD).Value = 44530
$ws.Cells.Item(190, $colJ).Value = 180
$ws.Cells.Item(190, $colK).Value = 5500
$ws.Cells.Item(190, $colL).Value = 5500
$ws.Cells.Item(190, $colM).Value = 5500
$ws.Cells.Item(190, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(190, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(190, $colP).Value = 1833
$ws.Cells.Item(190, $colQ).Value = 3

$ws.Cells.Item(191, $colD).Value = 44483
$ws.Cells.Item(191, $colJ).Value = 80
$ws.Cells.Item(191, $colK).Value = 4500
$ws.Cells.Item(191, $colL).Value = 4500
$ws.Cells.Item(191, $colM).Value = 4500
$ws.Cells.Item(191, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(191, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(191, $colP).Value = 1500
$ws.Cells.Item(191, $colQ).Value = 3

$ws.Cells.Item(192, $colD).Value = 44294
$ws.Cells.Item(192, $colJ).Value = 70
$ws.Cells.Item(192, $colK).Value = 5000
$ws.Cells.Item(192, $colL).Value = 5000
$ws.Cells.Item(192, $colM).Value = 5000
$ws.Cells.Item(192, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(192, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(192, $colP).Value = 1667
$ws.Cells.Item(192, $colQ).Value = 3

$ws.Cells.Item(193, $colD).Value = 44617
$ws.Cells.Item(193, $colJ).Value = 180
$ws.Cells.Item(193, $colK).Value = 5000
$ws.Cells.Item(193, $colL).Value = 5000
$ws.Cells.Item(193, $colM).Value = 5000
$ws.Cells.Item(193, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(193, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(193, $colP).Value = 1667
$ws.Cells.Item(193, $colQ).Value = 3

$ws.Cells.Item(194, $colD).Value = 44264
$ws.Cells.Item(194, $colJ).Value = 80
$ws.Cells.Item(194, $colK).Value = 4000
$ws.Cells.Item(194, $colL).Value = 4000
$ws.Cells.Item(194, $colM).Value = 4000
$ws.Cells.Item(194, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(194, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(194, $colP).Value = 2000
$ws.Cells.Item(194, $colQ).Value = 2

$ws.Cells.Item(195, $colD).Value = 44264
$ws.Cells.Item(195, $colJ).Value = 80
$ws.Cells.Item(195, $colK).Value = 5000
$ws.Cells.Item(195, $colL).Value = 5000
$ws.Cells.Item(195, $colM).Value = 5000
$ws.Cells.Item(195, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(195, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(195, $colP).Value = 1667
$ws.Cells.Item(195, $colQ).Value = 3

$ws.Cells.Item(196, $colD).Value = 44232
$ws.Cells.Item(196, $colJ).Value = 100
$ws.Cells.Item(196, $colK).Value = 5000
$ws.Cells.Item(196, $colL).Value = 5000
$ws.Cells.Item(196, $colM).Value = 5000
$ws.Cells.Item(196, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(196, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(196, $colP).Value = 2500
$ws.Cells.Item(196, $colQ).Value = 2

$ws.Cells.Item(197, $colD).Value = 44279
$ws.Cells.Item(197, $colJ).Value = 30
$ws.Cells.Item(197, $colK).Value = 4000
$ws.Cells.Item(197, $colL).Value = 4000
$ws.Cells.Item(197, $colM).Value = 4000
$ws.Cells.Item(197, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(197, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(197, $colP).Value = 2000
$ws.Cells.Item(197, $colQ).Value = 2

$ws.Cells.Item(198, $colD).Value = 44330
$ws.Cells.Item(198, $colJ).Value = 180
$ws.Cells.Item(198, $colK).Value = 4000
$ws.Cells.Item(198, $colL).Value = 4000
$ws.Cells.Item(198, $colM).Value = 4000
$ws.Cells.Item(198, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(198, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(198, $colP).Value = 1333
$ws.Cells.Item(198, $colQ).Value = 3

$ws.Cells.Item(199, $colD).Value = 44504
$ws.Cells.Item(199, $colJ).Value = 60
$ws.Cells.Item(199, $colK).Value = 5000
$ws.Cells.Item(199, $colL).Value = 5000
$ws.Cells.Item(199, $colM).Value = 5000
$ws.Cells.Item(199, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(199, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(199, $colP).Value = 1667
$ws.Cells.Item(199, $colQ).Value = 3

$ws.Cells.Item(200, $colD).Value = 44257
$ws.Cells.Item(200, $colJ).Value = 100
$ws.Cells.Item(200, $colK).Value = 4000
$ws.Cells.Item(200, $colL).Value = 4000
$ws.Cells.Item(200, $colM).Value = 4000
$ws.Cells.Item(200, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(200, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(200, $colP).Value = 2000
$ws.Cells.Item(200, $colQ).Value = 2

$ws.Cells.Item(201, $colD).Value = 44257
$ws.Cells.Item(201, $colJ).Value = 100
$ws.Cells.Item(201, $colK).Value = 5000
$ws.Cells.Item(201, $colL).Value = 5000
$ws.Cells.Item(201, $colM).Value = 5000
$ws.Cells.Item(201, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(201, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(201, $colP).Value = 1667
$ws.Cells.Item(201, $colQ).Value = 3

$ws.Cells.Item(202, $colD).Value = 44301
$ws.Cells.Item(202, $colJ).Value = 80
$ws.Cells.Item(202, $colK).Value = 4500
$ws.Cells.Item(202, $colL).Value = 4500
$ws.Cells.Item(202, $colM).Value = 4500
$ws.Cells.Item(202, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(202, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(202, $colP).Value = 1500
$ws.Cells.Item(202, $colQ).Value = 3

$ws.Cells.Item(203, $colD).Value = 44370
$ws.Cells.Item(203, $colJ).Value = 20
$ws.Cells.Item(203, $colK).Value = 4000
$ws.Cells.Item(203, $colL).Value = 4000
$ws.Cells.Item(203, $colM).Value = 4000
$ws.Cells.Item(203, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(203, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(203, $colP).Value = 1333
$ws.Cells.Item(203, $colQ).Value = 3

$ws.Cells.Item(204, $colD).Value = 44487
$ws.Cells.Item(204, $colJ).Value = 90
$ws.Cells.Item(204, $colK).Value = 5000
$ws.Cells.Item(204, $colL).Value = 5000
$ws.Cells.Item(204, $colM).Value = 5000
$ws.Cells.Item(204, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(204, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(204, $colP).Value = 1667
$ws.Cells.Item(204, $colQ).Value = 3

$ws.Cells.Item(205, $colD).Value = 44385
$ws.Cells.Item(205, $colJ).Value = 30
$ws.Cells.Item(205, $colK).Value = 4000
$ws.Cells.Item(205, $colL).Value = 4000
$ws.Cells.Item(205, $colM).Value = 4000
$ws.Cells.Item(205, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(205, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(205, $colP).Value = 1333
$ws.Cells.Item(205, $colQ).Value = 3

$ws.Cells.Item(206, $colD).Value = 44236
$ws.Cells.Item(206, $colJ).Value = 100
$ws.Cells.Item(206, $colK).Value = 5000
$ws.Cells.Item(206, $colL).Value = 5000
$ws.Cells.Item(206, $colM).Value = 5000
$ws.Cells.Item(206, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(206, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(206, $colP).Value = 2500
$ws.Cells.Item(206, $colQ).Value = 2

$ws.Cells.Item(207, $colD).Value = 44221
$ws.Cells.Item(207, $colJ).Value = 100
$ws.Cells.Item(207, $colK).Value = 4000
$ws.Cells.Item(207, $colL).Value = 4000
$ws.Cells.Item(207, $colM).Value = 4000
$ws.Cells.Item(207, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(207, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(207, $colP).Value = 2000
$ws.Cells.Item(207, $colQ).Value = 2

$ws.Cells.Item(208, $colD).Value = 44272
$ws.Cells.Item(208, $colJ).Value = 20
$ws.Cells.Item(208, $colK).Value = 4000
$ws.Cells.Item(208, $colL).Value = 4000
$ws.Cells.Item(208, $colM).Value = 4000
$ws.Cells.Item(208, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(208, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(208, $colP).Value = 2000
$ws.Cells.Item(208, $colQ).Value = 2

$ws.Cells.Item(209, $colD).Value = 44229
$ws.Cells.Item(209, $colJ).Value = 200
$ws.Cells.Item(209, $colK).Value = 4000
$ws.Cells.Item(209, $colL).Value = 5000
$ws.Cells.Item(209, $colM).Value = 4500
$ws.Cells.Item(209, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(209, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(209, $colP).Value = 2250
$ws.Cells.Item(209, $colQ).Value = 2

$ws.Cells.Item(210, $colD).Value = 44214
$ws.Cells.Item(210, $colJ).Value = 50
$ws.Cells.Item(210, $colK).Value = 4000
$ws.Cells.Item(210, $colL).Value = 4000
$ws.Cells.Item(210, $colM).Value = 4000
$ws.Cells.Item(210, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(210, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(210, $colP).Value = 2000
$ws.Cells.Item(210, $colQ).Value = 2

$ws.Cells.Item(211, $colD).Value = 44299
$ws.Cells.Item(211, $colJ).Value = 80
$ws.Cells.Item(211, $colK).Value = 5000
$ws.Cells.Item(211, $colL).Value = 5000
$ws.Cells.Item(211, $colM).Value = 5000
$ws.Cells.Item(211, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(211, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(211, $colP).Value = 2500
$ws.Cells.Item(211, $colQ).Value = 2

$ws.Cells.Item(212, $colD).Value = 44299
$ws.Cells.Item(212, $colJ).Value = 140
$ws.Cells.Item(212, $colK).Value = 4500
$ws.Cells.Item(212, $colL).Value = 4500
$ws.Cells.Item(212, $colM).Value = 4500
$ws.Cells.Item(212, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(212, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(212, $colP).Value = 1500
$ws.Cells.Item(212, $colQ).Value = 3

$ws.Cells.Item(213, $colD).Value = 44610
$ws.Cells.Item(213, $colJ).Value = 150
$ws.Cells.Item(213, $colK).Value = 5000
$ws.Cells.Item(213, $colL).Value = 5000
$ws.Cells.Item(213, $colM).Value = 5000
$ws.Cells.Item(213, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(213, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(213, $colP).Value = 1667
$ws.Cells.Item(213, $colQ).Value = 3

$ws.Cells.Item(214, $colD).Value = 44312
$ws.Cells.Item(214, $colJ).Value = 20
$ws.Cells.Item(214, $colK).Value = 5000
$ws.Cells.Item(214, $colL).Value = 5000
$ws.Cells.Item(214, $colM).Value = 5000
$ws.Cells.Item(214, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(214, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(214, $colP).Value = 2500
$ws.Cells.Item(214, $colQ).Value = 2

$ws.Cells.Item(215, $colD).Value = 44522
$ws.Cells.Item(215, $colJ).Value = 60
$ws.Cells.Item(215, $colK).Value = 5000
$ws.Cells.Item(215, $colL).Value = 5000
$ws.Cells.Item(215, $colM).Value = 5000
$ws.Cells.Item(215, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(215, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(215, $colP).Value = 1667
$ws.Cells.Item(215, $colQ).Value = 3

$ws.Cells.Item(216, $colD).Value = 44277
$ws.Cells.Item(216, $colJ).Value = 80
$ws.Cells.Item(216, $colK).Value = 4000
$ws.Cells.Item(216, $colL).Value = 4000
$ws.Cells.Item(216, $colM).Value = 4000
$ws.Cells.Item(216, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(216, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(216, $colP).Value = 2000
$ws.Cells.Item(216, $colQ).Value = 2

$ws.Cells.Item(217, $colD).Value = 44258
$ws.Cells.Item(217, $colJ).Value = 20
$ws.Cells.Item(217, $colK).Value = 4000
$ws.Cells.Item(217, $colL).Value = 4000
$ws.Cells.Item(217, $colM).Value = 4000
$ws.Cells.Item(217, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(217, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(217, $colP).Value = 2000
$ws.Cells.Item(217, $colQ).Value = 2

$ws.Cells.Item(218, $colD).Value = 44390
$ws.Cells.Item(218, $colJ).Value = 180
$ws.Cells.Item(218, $colK).Value = 4500
$ws.Cells.Item(218, $colL).Value = 4500
$ws.Cells.Item(218, $colM).Value = 4500
$ws.Cells.Item(218, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(218, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(218, $colP).Value = 1500
$ws.Cells.Item(218, $colQ).Value = 3

$ws.Cells.Item(219, $colD).Value = 44349
$ws.Cells.Item(219, $colJ).Value = 30
$ws.Cells.Item(219, $colK).Value = 3500
$ws.Cells.Item(219, $colL).Value = 3500
$ws.Cells.Item(219, $colM).Value = 3500
$ws.Cells.Item(219, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(219, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(219, $colP).Value = 1167
$ws.Cells.Item(219, $colQ).Value = 3

$ws.Cells.Item(220, $colD).Value = 44285
$ws.Cells.Item(220, $colJ).Value = 100
$ws.Cells.Item(220, $colK).Value = 4000
$ws.Cells.Item(220, $colL).Value = 4000
$ws.Cells.Item(220, $colM).Value = 4000
$ws.Cells.Item(220, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(220, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(220, $colP).Value = 2000
$ws.Cells.Item(220, $colQ).Value = 2

$ws.Cells.Item(221, $colD).Value = 44285
$ws.Cells.Item(221, $colJ).Value = 150
$ws.Cells.Item(221, $colK).Value = 5000
$ws.Cells.Item(221, $colL).Value = 5000
$ws.Cells.Item(221, $colM).Value = 5000
$ws.Cells.Item(221, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(221, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(221, $colP).Value = 1667
$ws.Cells.Item(221, $colQ).Value = 3

$ws.Cells.Item(222, $colD).Value = 44498
$ws.Cells.Item(222, $colJ).Value = 160
$ws.Cells.Item(222, $colK).Value = 4500
$ws.Cells.Item(222, $colL).Value = 4500
$ws.Cells.Item(222, $colM).Value = 4500
$ws.Cells.Item(222, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(222, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(222, $colP).Value = 1500
$ws.Cells.Item(222, $colQ).Value = 3

$ws.Cells.Item(223, $colD).Value = 44418
$ws.Cells.Item(223, $colJ).Value = 180
$ws.Cells.Item(223, $colK).Value = 5000
$ws.Cells.Item(223, $colL).Value = 5000
$ws.Cells.Item(223, $colM).Value = 5000
$ws.Cells.Item(223, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(223, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(223, $colP).Value = 1667
$ws.Cells.Item(223, $colQ).Value = 3

$ws.Cells.Item(224, $colD).Value = 44595
$ws.Cells.Item(224, $colJ).Value = 80
$ws.Cells.Item(224, $colK).Value = 5000
$ws.Cells.Item(224, $colL).Value = 5000
$ws.Cells.Item(224, $colM).Value = 5000
$ws.Cells.Item(224, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(224, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(224, $colP).Value = 1667
$ws.Cells.Item(224, $colQ).Value = 3

$ws.Cells.Item(225, $colD).Value = 44628
$ws.Cells.Item(225, $colJ).Value = 180
$ws.Cells.Item(225, $colK).Value = 6000
$ws.Cells.Item(225, $colL).Value = 6000
$ws.Cells.Item(225, $colM).Value = 6000
$ws.Cells.Item(225, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(225, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(225, $colP).Value = 2000
$ws.Cells.Item(225, $colQ).Value = 3

$ws.Cells.Item(226, $colD).Value = 44552
$ws.Cells.Item(226, $colJ).Value = 20
$ws.Cells.Item(226, $colK).Value = 7000
$ws.Cells.Item(226, $colL).Value = 7000
$ws.Cells.Item(226, $colM).Value = 7000
$ws.Cells.Item(226, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(226, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(226, $colP).Value = 3500
$ws.Cells.Item(226, $colQ).Value = 2

$ws.Cells.Item(227, $colD).Value = 44544
$ws.Cells.Item(227, $colJ).Value = 80
$ws.Cells.Item(227, $colK).Value = 6000
$ws.Cells.Item(227, $colL).Value = 6000
$ws.Cells.Item(227, $colM).Value = 6000
$ws.Cells.Item(227, $colN).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(227, $colO).Value = 'Región de La Araucanía'
$ws.Cells.Item(227, $colP).Value = 3000
$ws.Cells.Item(227, $colQ).Value = 2

$ws.Cells.Item(228, $colD).Value = 44544
$ws.Cells.Item(228, $colJ).Value = 160
$ws.Cells.Item(228, $colK).Value = 5000
$ws.Cells.Item(228, $colL).Value = 5000
$ws.Cells.Item(228, $colM).Value = 5000
$ws.Cells.Item(228, $colN).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(228, $colO).Value = 'Región Metropolitana'
$ws.Cells.Item(228, $colP).Value = 1667
$ws.Cells.Item(228, $colQ).Value = 3

# Rows 227-228 are brand new rows appended at the bottom of the table; populate the
# columns that stay constant across this entire subset (A, B, C, E, F, G, H, I, R),
# and apply the same date number format used by column D elsewhere in the sheet.
foreach ($r in 227,228) {
    $ws.Cells.Item($r, 1).Value = 4
    $ws.Cells.Item($r, 2).Value = 'Feria Lagunitas de Puerto Montt'
    $ws.Cells.Item($r, 3).Value = 'Los Lagos'
    $ws.Cells.Item($r, $colD).NumberFormat = $ws.Cells.Item(226, $colD).NumberFormat
    $ws.Cells.Item($r, 5).Value = 10
    $ws.Cells.Item($r, 6).Value = 100112044
    $ws.Cells.Item($r, 7).Value = 'Perejil'
    $ws.Cells.Item($r, 8).Value = 'Sin especificar'
    $ws.Cells.Item($r, 9).Value = 'Primera'
    $ws.Cells.Item($r, 18).Value = 'Hortaliza'
}

